$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.321.84'
$ws.Range("E2").Value = '  -2.96%  '
$ws.Range("D3").Value = '1.934.84'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '250.27'
$ws.Range("E5").Value = '  -1.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7217'
$ws.Range("E6").Value = '  -7.35%  '
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3301'
$ws.Range("E8").Value = '  -5.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '27.72'
$ws.Range("E9").Value = '  -0.97%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07188'
$ws.Range("E10").Value = '  +1.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8086'
$ws.Range("E11").Value = '  -3.89%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08090'
$ws.Range("E12").Value = '  -1.18%  '
$ws.Range("D13").Value = '1.937.21'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.485'
$ws.Range("E14").Value = '  -2.67%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '94.35'
$ws.Range("E15").Value = '  -6.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.17'
$ws.Range("E16").Value = '  -0.58%  '
$ws.Range("D17").Value = '30.345.57'
$ws.Range("E17").Value = '  -2.89%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008254'
$ws.Range("E18").Value = '  +3.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '249.73'
$ws.Range("E19").Value = '  -8.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.895'
$ws.Range("E20").Value = '  -1.49%  '
$ws.Range("D21").Value = '2.192.46'
$ws.Range("E21").Value = '  -2.81%  '
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.974'
$ws.Range("E24").Value = '  -1.69%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.737'
$ws.Range("E25").Value = '  -3.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.53'
$ws.Range("E26").Value = '  -0.66%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.372'
$ws.Range("E27").Value = '  -1.38%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.22'
$ws.Range("E28").Value = '  -3.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1317'
$ws.Range("E29").Value = '  -7.52%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.565'
$ws.Range("E30").Value = '  -1.67%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.353'
$ws.Range("E31").Value = '  -1.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.425'
$ws.Range("E32").Value = '  -4.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.167'
$ws.Range("E33").Value = '  -6.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05187'
$ws.Range("E34").Value = '  -2.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.278'
$ws.Range("E35").Value = '  +2.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7475'
$ws.Range("E36").Value = '  -5.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.743'
$ws.Range("E37").Value = '  -0.89%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01979'
$ws.Range("E38").Value = '  -1.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.825'
$ws.Range("E39").Value = '  -3.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '79.42'
$ws.Range("E40").Value = '  -5.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.369'
$ws.Range("E41").Value = '  -5.91%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4530'
$ws.Range("E42").Value = '  -3.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.021'
$ws.Range("E43").Value = '  -5.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8466'
$ws.Range("E44").Value = '  -1.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  +0.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.64'
$ws.Range("E46").Value = '  -3.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.760'
$ws.Range("E47").Value = '  -2.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.448'
$ws.Range("E48").Value = '  -4.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.71'
$ws.Range("E49").Value = '  -2.77%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4177'
$ws.Range("E50").Value = '  -3.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06037'
$ws.Range("E51").Value = '  +0.28%  '
